$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting rows 11-112 down to 12-113
$ws.Rows.Item(11).Insert()

# Fill the new row 11 with data
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44503
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100101
$ws.Cells.Item(11, 8).Value = "Berries"
$ws.Cells.Item(11, 9).Value = 100101001
$ws.Cells.Item(11, 10).Value = "Arándano (blue)"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 350
$ws.Cells.Item(11, 14).Value = 13000
$ws.Cells.Item(11, 15).Value = 14000
$ws.Cells.Item(11, 16).Value = 13429
$ws.Cells.Item(11, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 6714
$ws.Cells.Item(11, 20).Value = 2
